$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from H1 to I1:J1 (bold header style), then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 11
$ws.Range("J5").Value = 11

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9
